$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1/24
$ws.Range("C3").Value = 1.25/24
$ws.Range("C4").Value = 1.25/24
$ws.Range("C5").Value = 1.5/24

$ws.Range("C3").Select()
